$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell contents ---

# Row 7: E7 cleared (task moved elsewhere)
$ws.Range("E7").ClearContents()

# New G column entry for row 4 (re-uses existing shared string)
$ws.Range("G4").Value = "Obter as Code Base Metrics (José Pereira)"

# New G column entry for row 5
$ws.Range("G5").Value = "Identificar os Code Smells (José Pereira)"

# Row 8: B8 text change
$ws.Range("B8").Value = "Identificar 3 GOF Design Patterns"

# Row 9: B9 text change
$ws.Range("B9").Value = "Identificar 3 Code Smells"

# Row 10 stays "Organizar o repositório no Github" (unchanged)

# New row 11: reinsert the "Reunir o grupo..." task that moved out of B9
$ws.Range("B11").Value = "Reunir o grupo para organizar a entrega da 2ª fase do projeto"

# Row 6: D6 / E6 content change
$ws.Range("D6").Value = "Identificar 3 Design Patterns (Todos têm de fazer)"
$ws.Range("E6").Value = "Reunir o grupo para organizar a entrega da 2ª fase do projeto (José Pereira)"

# --- Column widths (D and E got resized) ---
$ws.Columns.Item(4).ColumnWidth = 56.17
$ws.Columns.Item(5).ColumnWidth = 66.67

# --- Selection / active cell ---
$ws.Range("D6").Select() | Out-Null
